$wb = $excel.ActiveWorkbook
$loginSheet = $wb.Worksheets.Item("LoginTest")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$newSheet.Name = "RegisterTest"
$newSheet.Range("B1").Value = "LastName"
$newSheet.Range("A1").Value = "FirstName"
$newSheet.Range("A2").Value = "Test"
$newSheet.Range("B2").Value = "Here"
$newSheet.Rows(1).Font.Bold = $true
$newSheet.Rows(1).Interior.ColorIndex = -4142
